# "Added a good number of equip. to database. Also added some temp techs
#  to Austria to avoid so many errors"
#
# The change adds two more generations of the Light Strike Fighter
# equipment/tech pair (L_Strike_fighter_equipment_2/3 + L_Strike_fighter2/3)
# right below the existing "_1" row, pushing the Helicopter block that used
# to start at row 182 down by three rows (to row 185) without altering it.
# It also removes the bold emphasis that had been applied (seemingly by
# mistake) to the "NB! All codes are case-sensitive..." note in row 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 ("NB! All codes are case-sensitive...") should not be bold.
$ws.Range("A4").Font.Bold = $false

# Make room for the two new Strike Fighter rows: push everything from the
# old row 182 ("early_helicopter" / Transport Helicopter row) down by 3.
$ws.Rows.Item(182).Resize(3).Insert()

# Row 179 already has "L_Strike_fighter_equipment_1" in column A - add its
# matching tech code + generation year.
$ws.Range("B179").Value = "L_Strike_fighter1"
$ws.Range("C179").Value = 1975

# Row 180 already has "L_Strike_fighter_equipment_2" in column A - add its
# matching tech code + generation year.
$ws.Range("B180").Value = "L_Strike_fighter2"
$ws.Range("C180").Value = 1995

# Row 181 is brand new - the 3rd generation Light Strike Fighter equipment.
$ws.Range("A181").Value = "L_Strike_fighter_equipment_3"
$ws.Range("B181").Value = "L_Strike_fighter3"
$ws.Range("C181").Value = 2015
